$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 ("Experimental" value): true -> false.
# A bare Value assignment of "false" gets auto-typed as a boolean by the
# engine (like real Excel would for an unformatted cell), so force a text
# literal with a leading quote-prefix, then restore the original cell
# formatting (which the quote-prefix bumps to a new style) by pasting the
# formats from the untouched cell directly above it.
$ws.Range("B7").Value = "'false"
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B8 ("Date" value): timestamp refresh - plain text, no coercion risk.
$ws.Range("B8").Value = "2025-10-03T16:37:46+01:00"
